$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update version numbers / labels that remain in place ---
# LZ4 version 1.8.2 -> 1.9.3
$ws.Range("C2").Value = "1.9.3"
# OpenSSL version 1.1.0h -> 1.1.1k
$ws.Range("C6").Value = "1.1.1k"

# --- Remove obsolete rows ---
# Remove the stb_image row (row 7)
$ws.Rows.Item(7).Delete()
# Remove the gsl row (was row 11, now row 10 after the previous delete)
$ws.Rows.Item(10).Delete()
# Remove the "*Source was taken from a post 2.3.2 release..." footnote row
# (was row 16, now row 14 after the previous two deletes)
$ws.Rows.Item(14).Delete()

# --- Update the vcredist_x64 entry (now at row 12) ---
$ws.Range("B12").Value = "Visual Studio 2019 (15.4.2) Runtimes"
$ws.Range("C12").Value = "14.28.29914.0"

# --- Fix up hyperlinks: row deletion above does not re-anchor hyperlink refs,
#     so rebuild the hyperlinks collection from scratch at the correct cells.
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://lz4.github.io/lz4/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "http://qt-project.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "http://qt-project.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "http://sqlite.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.openssl.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "http://rapidjson.org/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://go.microsoft.com/fwlink/?LinkId=746572") | Out-Null

# --- Selection as left by the edit ---
$ws.Range("C6").Select()
